$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("1er Parcial")
$ws.Range("C5").Value = 21
$ws.Range("D5").Value = 26.25
$ws.Range("E5").Value = 59
$ws.Range("F5").Value = 73.75
$ws.Range("I5").Value = 8.5
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("G9").Value = 16
$ws.Range("H9").Value = 16.49
$ws.Range("I9").Value = 7.6
$ws.Range("C18").Value = 47
$ws.Range("D18").Value = 19.03
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 80.97
$ws.Range("C22").Value = 21
$ws.Range("D22").Value = 10.29
$ws.Range("E22").Value = 135
$ws.Range("F22").Value = 66.18000000000001
$ws.Range("G22").Value = 48
$ws.Range("H22").Value = 23.53
$ws.Range("C23").Value = 44
$ws.Range("D23").Value = 38.94
$ws.Range("E23").Value = 69
$ws.Range("F23").Value = 61.06
$ws.Range("I23").Value = 8.1
$ws.Range("C24").Value = 3
$ws.Range("D24").Value = 3.85
$ws.Range("G24").Value = 16
$ws.Range("H24").Value = 20.51
$ws.Range("I24").Value = 6.4
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 180
$ws.Range("F25").Value = 100
$ws.Range("I25").Value = 8.1
$ws.Range("C27").Value = 32
$ws.Range("D27").Value = 27.59
$ws.Range("E27").Value = 84
$ws.Range("F27").Value = 72.41
$ws.Range("I27").Value = 7.6
$ws.Range("C30").Value = 52
$ws.Range("D30").Value = 24.76
$ws.Range("E30").Value = 148
$ws.Range("F30").Value = 70.48
$ws.Range("C39").Value = 13
$ws.Range("D39").Value = 9.49
$ws.Range("G39").Value = 50
$ws.Range("H39").Value = 36.5
$ws.Range("I39").Value = 6.2
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = 0
$ws.Range("G41").Value = 23
$ws.Range("H41").Value = 20
$ws.Range("I41").Value = 6.8
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 94
$ws.Range("F43").Value = 100
$ws.Range("I43").Value = 6.5
$ws.Range("C47").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 137
$ws.Range("F47").Value = 69.54000000000001
$ws.Range("G47").Value = 60
$ws.Range("H47").Value = 30.46
$ws.Range("I47").Value = 7
$ws.Range("C48").Value = 29
$ws.Range("D48").Value = 16.48
$ws.Range("E48").Value = 147
$ws.Range("F48").Value = 83.52

$ws = $wb.Worksheets.Item("2o Parcial")
$ws.Range("C4").Value = 64
$ws.Range("D4").Value = 39.02
$ws.Range("E4").Value = 73
$ws.Range("F4").Value = 44.51
$ws.Range("G4").Value = 91
$ws.Range("H4").Value = 55.49
$ws.Range("I4").Value = 7.6
$ws.Range("C5").Value = 24
$ws.Range("D5").Value = 30
$ws.Range("E5").Value = 56
$ws.Range("F5").Value = 70
$ws.Range("G5").Value = 3
$ws.Range("H5").Value = 3.75
$ws.Range("I5").Value = 7.8
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("G9").Value = 25
$ws.Range("H9").Value = 25.77
$ws.Range("C11").Value = 16
$ws.Range("D11").Value = 41.03
$ws.Range("E11").Value = 23
$ws.Range("F11").Value = 58.97
$ws.Range("G11").Value = 6
$ws.Range("H11").Value = 15.38
$ws.Range("I11").Value = 7
$ws.Range("C18").Value = 69
$ws.Range("D18").Value = 27.94
$ws.Range("E18").Value = 178
$ws.Range("F18").Value = 72.06
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = 8.91
$ws.Range("I18").Value = 7.7
$ws.Range("C22").Value = 65
$ws.Range("D22").Value = 31.86
$ws.Range("E22").Value = 131
$ws.Range("F22").Value = 64.22
$ws.Range("G22").Value = 52
$ws.Range("H22").Value = 25.49
$ws.Range("I22").Value = 8.300000000000001
$ws.Range("C23").Value = 68
$ws.Range("D23").Value = 60.18
$ws.Range("E23").Value = 45
$ws.Range("F23").Value = 39.82
$ws.Range("G23").Value = 24
$ws.Range("H23").Value = 21.24
$ws.Range("I23").Value = 8.4
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 30.77
$ws.Range("E24").Value = 54
$ws.Range("F24").Value = 69.23
$ws.Range("G24").Value = 21
$ws.Range("H24").Value = 26.92
$ws.Range("I24").Value = 7
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 180
$ws.Range("F25").Value = 100
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 8.199999999999999
$ws.Range("C26").Value = 16
$ws.Range("D26").Value = 17.78
$ws.Range("E26").Value = 70
$ws.Range("F26").Value = 77.78
$ws.Range("G26").Value = 20
$ws.Range("H26").Value = 22.22
$ws.Range("I26").Value = 8
$ws.Range("C27").Value = 51
$ws.Range("D27").Value = 43.97
$ws.Range("E27").Value = 65
$ws.Range("F27").Value = 56.03
$ws.Range("G27").Value = 19
$ws.Range("H27").Value = 16.38
$ws.Range("I27").Value = 8.199999999999999
$ws.Range("C30").Value = 78
$ws.Range("D30").Value = 37.14
$ws.Range("E30").Value = 127
$ws.Range("F30").Value = 60.48
$ws.Range("G30").Value = 31
$ws.Range("H30").Value = 14.76
$ws.Range("I30").Value = 7.7
$ws.Range("C39").Value = 70
$ws.Range("D39").Value = 51.09
$ws.Range("E39").Value = 50
$ws.Range("F39").Value = 36.5
$ws.Range("G39").Value = 74
$ws.Range("H39").Value = 54.01
$ws.Range("I39").Value = 6.6
$ws.Range("C41").Value = 28
$ws.Range("D41").Value = 24.35
$ws.Range("E41").Value = 64
$ws.Range("F41").Value = 55.65
$ws.Range("G41").Value = 51
$ws.Range("H41").Value = 44.35
$ws.Range("I41").Value = 7.3
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 94
$ws.Range("F43").Value = 100
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 6.5
$ws.Range("C47").Value = 87
$ws.Range("D47").Value = 44.16
$ws.Range("E47").Value = 108
$ws.Range("F47").Value = 54.82
$ws.Range("G47").Value = 89
$ws.Range("H47").Value = 45.18
$ws.Range("I47").Value = 8
$ws.Range("C48").Value = 29
$ws.Range("D48").Value = 16.48
$ws.Range("E48").Value = 147
$ws.Range("F48").Value = 83.52
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 7.4
$ws.Range("C49").Value = 68
$ws.Range("D49").Value = 55.28
$ws.Range("E49").Value = 55
$ws.Range("F49").Value = 44.72
$ws.Range("G49").Value = 68
$ws.Range("H49").Value = 55.28

$ws = $wb.Worksheets.Item("Final")
$ws.Range("E4").Value = 115
$ws.Range("F4").Value = 70.12
$ws.Range("G4").Value = 49
$ws.Range("H4").Value = 29.88
$ws.Range("C5").Value = 21
$ws.Range("D5").Value = 26.25
$ws.Range("E5").Value = 59
$ws.Range("F5").Value = 73.75
$ws.Range("I5").Value = 8.5
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("G9").Value = 25
$ws.Range("H9").Value = 25.77
$ws.Range("I9").Value = 7.6
$ws.Range("C18").Value = 47
$ws.Range("D18").Value = 19.03
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 80.97
$ws.Range("C22").Value = 21
$ws.Range("D22").Value = 10.29
$ws.Range("E22").Value = 146
$ws.Range("F22").Value = 71.56999999999999
$ws.Range("G22").Value = 37
$ws.Range("H22").Value = 18.14
$ws.Range("C23").Value = 44
$ws.Range("D23").Value = 38.94
$ws.Range("E23").Value = 69
$ws.Range("F23").Value = 61.06
$ws.Range("I23").Value = 8.1
$ws.Range("C24").Value = 3
$ws.Range("D24").Value = 3.85
$ws.Range("G24").Value = 13
$ws.Range("H24").Value = 16.67
$ws.Range("I24").Value = 6.4
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 180
$ws.Range("F25").Value = 100
$ws.Range("I25").Value = 8.1
$ws.Range("C27").Value = 32
$ws.Range("D27").Value = 27.59
$ws.Range("E27").Value = 84
$ws.Range("F27").Value = 72.41
$ws.Range("I27").Value = 7.6
$ws.Range("C30").Value = 52
$ws.Range("D30").Value = 24.76
$ws.Range("E30").Value = 147
$ws.Range("F30").Value = 70
$ws.Range("C39").Value = 13
$ws.Range("D39").Value = 9.49
$ws.Range("G39").Value = 56
$ws.Range("H39").Value = 40.88
$ws.Range("I39").Value = 6.2
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = 0
$ws.Range("E41").Value = 77
$ws.Range("F41").Value = 66.95999999999999
$ws.Range("G41").Value = 38
$ws.Range("H41").Value = 33.04
$ws.Range("I41").Value = 6.8
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 94
$ws.Range("F43").Value = 100
$ws.Range("I43").Value = 6.5
$ws.Range("C47").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 139
$ws.Range("F47").Value = 70.56
$ws.Range("G47").Value = 58
$ws.Range("H47").Value = 29.44
$ws.Range("I47").Value = 7
$ws.Range("C48").Value = 29
$ws.Range("D48").Value = 16.48
$ws.Range("E48").Value = 147
$ws.Range("F48").Value = 83.52
